$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: fix the sentence "...the allocations of marks shows:" ->
#         "...the allocations of marks shows as below:" and drop the
#         (now-resolved) grammar proofing marks around "shows".
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.Contains("allocations of marks")) {
        $target = $cand
        break
    }
}

if ($target -ne $null) {
    $pXml = $target.Range.WordOpenXML
    $openStart = $pXml.IndexOf("<w:p ")
    if ($openStart -eq -1) {
        $openStart = $pXml.IndexOf("<w:p>")
    }
    $openEnd = $pXml.IndexOf(">", $openStart) + 1
    $pOpenTag = $pXml.Substring($openStart, $openEnd - $openStart)

    $newParaXml = $pOpenTag + `
        '<w:r><w:t>As a team, we have all voted to equally share all the marks between us 4 members. With this being said, the allocations of marks show</w:t></w:r>' + `
        '<w:r><w:t>s as below</w:t></w:r>' + `
        '<w:r><w:t>:</w:t></w:r>' + `
        '</w:p>'

    $target.Range.InsertXML($newParaXml)
}

# ---------------------------------------------------------------------------
# Part 2: the three watermark pictures living in the section headers need an
# explicit <w10:wrap anchorx="margin" anchory="margin"/> added inside their
# VML <v:shape> so the watermark stays anchored/centred on the page margin.
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists -eq $false) { continue }
    if ($hdr.Shapes.Count -eq 0) { continue }

    $hXml = $hdr.Range.WordOpenXML
    $pStart = $hXml.IndexOf("<w:p ")
    if ($pStart -eq -1) {
        $pStart = $hXml.IndexOf("<w:p>")
    }
    $pEnd = $hXml.IndexOf("</w:p>", $pStart) + "</w:p>".Length
    $pFrag = $hXml.Substring($pStart, $pEnd - $pStart)

    if ($pFrag.Contains("</v:shape>") -and (-not $pFrag.Contains("w10:wrap"))) {
        $newFrag = $pFrag.Replace("</v:shape>", '<w10:wrap anchorx="margin" anchory="margin"/></v:shape>')
        $hdr.Range.InsertXML($newFrag)
    }
}

Write-Output "Done"
